$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50 (pushes former rows 50..84 down to 51..85),
# matching the weekly price-entry that was added to the source table.
$ws.Rows.Item(50).Insert()

$ws.Range("A50").Value = 5
$ws.Range("B50").Value = 'Macroferia Regional de Talca'
$ws.Range("C50").Value = 'Maule'
$ws.Range("D50").Value = 44879
$ws.Range("E50").Value = 7
$ws.Range("F50").Value = 300000000
$ws.Range("G50").Value = 'Espárragos'
$ws.Range("H50").Value = 'Sin especificar'
$ws.Range("I50").Value = 'Primera'
$ws.Range("J50").Value = 3000
$ws.Range("K50").Value = 1000
$ws.Range("L50").Value = 1000
$ws.Range("M50").Value = 1000
$ws.Range("N50").Value = '$/kilo'
$ws.Range("O50").Value = 'Provincia de Linares'
$ws.Range("P50").Value = 1000
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = 'Hortaliza'
